$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    2.0714786385807313,
    2.08376792163845,
    2.0947616571668086,
    2.1045003371691613,
    2.1130244536488636,
    2.1203744986092743,
    2.1265909640537424,
    2.1317143419856275,
    2.1357851244082857,
    2.138843803325065,
    2.1409308707393295,
    2.14208681865443,
    2.1423521390737212,
    2.1417673240005604,
    2.1403728654383025,
    2.1382092553903003,
    2.1353169858599115,
    2.13173654885049,
    2.1275084363653916,
    2.122673140407972,
    2.1172711529815866,
    2.111342966089587,
    2.1049290717353335,
    2.09806996192218,
    2.090806128653478,
    2.0831780639325865,
    2.075226259762861,
    2.068035910300791,
    2.060355485788047,
    2.0522128252047116,
    2.04355579521056,
    2.0345685035340675,
    2.025352153715997,
    2.0162404084975365,
    2.007185319252745,
    1.9980247107622389,
    1.9890940124067025,
    1.9805405929009918,
    1.9720188973764308,
    1.9638533022628248,
    1.955895489624356,
    1.9480981889331677,
    1.9405164557383392,
    1.9332093033364919,
    1.9262380182598737,
    1.9198069165174698,
    1.913798793890147,
    1.908485627111801,
    1.9039477220104988,
    1.9000768187451968,
    1.8960149827229678,
    1.8919735859592552,
    1.8877732252343564,
    1.8833558588414552,
    1.879131289615586,
    1.8748554710865024,
    1.8709077004169816,
    1.8672554650757747,
    1.863881282378746,
    1.8609593447277084,
    1.8581740453428286,
    1.8554672035442574,
    1.8529888626868778,
    1.8505925490313315,
    1.8482452726868992,
    1.8459198278325866,
    1.8435147658079623,
    1.8411558605610663,
    1.838893851534868,
    1.8367696670764135,
    1.8346058784935613,
    1.8322924395334634,
    1.8295478719278915,
    1.8265411981902788,
    1.8232017259503273,
    1.8199299528975255,
    1.8166841486844145,
    1.8131577211598955,
    1.8094952022754058,
    1.8059796856038153,
    1.802737332161435,
    1.7995680291209464,
    1.796421620222504,
    1.7933252035564768,
    1.790262900284055,
    1.787248747512872,
    1.7842430447143351,
    1.7811895809017306,
    1.778130870442744,
    1.7751525660793503,
    1.7720877359253169,
    1.7691192834196183,
    1.7659863387874382,
    1.7630007686017597,
    1.7602698166725672,
    1.7578822152606595,
    1.7557075967503168,
    1.7538183335540127,
    1.7516404027946817,
    1.749546473456707,
    1.7478817259871327,
    1.746410620950618,
    1.7453105205383788,
    1.7443920195829108,
    1.7433732768837813,
    1.741988511867396,
    1.7404823711261532,
    1.738737214009909,
    1.7370009197180114,
    1.7350919593752103,
    1.7329621809457552,
    1.730428813084738,
    1.7273087995816532,
    1.72354844390336,
    1.7195137098010163,
    1.7152509642988079,
    1.7109685457180255,
    1.7067972304883725,
    1.7027479924436728,
    1.6987630441282842,
    1.6946343165602906,
    1.6901722386471552,
    1.685209588537914,
    1.6798904289673897,
    1.6747492646890567,
    1.6697523410942534,
    1.6651564928218152,
    1.66118160273452,
    1.6575789257056532,
    1.6548351294178216,
    1.6525159858524614,
    1.6505567841382274,
    1.6489301151182876,
    1.6472264513261434,
    1.645626219171713,
    1.644023133765001,
    1.642280422621648,
    1.640518069317899,
    1.6387362923477358,
    1.6371617173365305,
    1.6357941255916202,
    1.6349236502131657,
    1.634536364387115,
    1.6345035701189388,
    1.634825350998438,
    1.6354069655149515,
    1.6361800475183614,
    1.6369927415989123,
    1.6376335398478197,
    1.6381207408647334,
    1.63841931662294,
    1.6390754171675983,
    1.6398996618829145,
    1.6403674714201517,
    1.6407942177122505,
    1.6407925863800799,
    1.6405496201048069,
    1.6401301611826167,
    1.6396563164627542,
    1.6392517516444092,
    1.638826775732827,
    1.6383471873978916,
    1.6379307585121452,
    1.637418499310112,
    1.6369213022637252,
    1.6366155378380525,
    1.63668628099072,
    1.637096734788431,
    1.6378940587428534,
    1.6386795603534983,
    1.6392246589593877,
    1.6395421467863336,
    1.6396570109545938,
    1.6396840628549099,
    1.641043158151257,
    1.6424441897345579,
    1.6438715367321792,
    1.6453095782714908,
    1.6467426934798572,
    1.64815526148465,
    1.6495316614132354,
    1.6508562723929816,
    1.6521134735512566,
    1.6532876440154294,
    1.654363162912866,
    1.6553244093709354,
    1.656155762517006,
    1.656841601478445,
    1.657366305382621,
    1.6577142533569018,
    1.6578698245286554,
    1.6578173980252493,
    1.6575413529740517,
    1.6570260685024312,
    1.6562559237377548,
    1.6552152978073913,
    1.6538885698387085,
    1.6522601189590738,
    1.6503143242958556,
    1.648035564976422
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
